$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Portfolio")

$data = @(
    @("394670", "TIGER 글로벌리튬&2차전지SOLACTIVE(합성)", "글로벌주식", 85.72653, 0.05),
    @("292150", "TIGER 코리아TOP10", "국내주식", 40.06132, 0.0588),
    @("483420", "ACE 애플밸류체인액티브", "미국주식", 298.8439, 0.2585),
    @("245350", "TIGER 유로스탁스배당30", "유럽주식", 213.8306, 0.2064),
    @("469160", "ACE 일본반도체", "일본주식", 71.72225, 0.0785),
    @("220130", "SOL 차이나강소기업CSI500(합성 H)", "중국본토주식", 97.79045000000001, 0.0343),
    @("419430", "KODEX 차이나2차전지MSCI(합성)", "중국역외주식", 59.7561, 0.0343),
    @("105010", "TIGER 라틴35", "신흥국주식", 329.6137, 0.1152),
    @("455960", "RISE 미국달러SOFR금리액티브(합성)", "선진국채권1", 265.0505, 0.164)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]

    # Column A holds numeric-looking ticker codes that must stay text
    # (e.g. leading zeros must be preserved). Force text storage, write
    # the value, then reset the style back to Normal so no stray
    # number-format/style is left behind on the cell.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $vals[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}
